$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.545.88"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "'3.392.45"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'403.74"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'130.30"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.590"
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.681"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "'0.129"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "'41.48"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'8.32"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").Value = "'3.393.36"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "'11.59"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "'61.491.65"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "'0.0000142"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'3.15"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").Value = "'82.99"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'310.62"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'12.68"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +5.72%  "
$ws.Range("D26").Value = "'8.04"
$ws.Range("E26").Value = "  +6.36%  "
$ws.Range("D27").Value = "'29.34"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'8.02"
$ws.Range("E28").Value = "  -8.27%  "
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").Value = "'43.79"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "'11.25"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").Value = "'51.22"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'2.96"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Value = "'3.33"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "'0.311"
$ws.Range("E40").Value = "  +8.91%  "
$ws.Range("D41").Value = "'139.91"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").Value = "'1.94"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "'3.91"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'16.63"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'21.03"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").Value = "'2.093.15"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = "  +11.95%  "
